$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6) - shifts GENES..REFERENCES(DOI)
# header block (old F:L) one column to the right (new G:M), and copies the
# formatting (style + column width) of the column to the left (E) onto it,
# matching Excel's default "Insert" behaviour.
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# New column header
$ws.Range("F1").Value = "CELLS"

# The CT (cell-type) lists that used to live in the CHILDREN column (E) for
# the AS rows 5-8 now belong under the new CELLS column (F).
$ws.Range("F5").Value = $ws.Range("E5").Value2
$ws.Range("E5").ClearContents()

$ws.Range("F6").Value = $ws.Range("E6").Value2
$ws.Range("E6").ClearContents()

$ws.Range("F7").Value = $ws.Range("E7").Value2
$ws.Range("E7").ClearContents()

$ws.Range("F8").Value = $ws.Range("E8").Value2
$ws.Range("E8").ClearContents()

# Match the final selection recorded in the saved workbook.
$ws.Range("F1").Select()
